# Update the "dSF" column (F) values for several rows to match the
# repulled / recalculated data, per commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = 1
$ws.Range("F32").Value = -1
$ws.Range("F39").Value = 1
$ws.Range("F41").Value = 0
$ws.Range("F50").Value = -2
$ws.Range("F51").Value = -1
$ws.Range("F58").Value = 0
$ws.Range("F59").Value = -10
$ws.Range("F62").Value = -1
